# feat: add 2022-Q4 data
#
# Before: sheet "总计" (totals) + sheet "2022-Q3" (fund holdings for Q3)
# After : sheet "总计" (totals, now with a Q4 row) + sheet "2022-Q4" (new fund
#         holdings data) + sheet "2022-Q3" (the original, unchanged, Q3 data)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: update the "总计" (totals) sheet - insert the new 2022-Q4 row above
# the existing 2022-Q3 row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Move the existing 2022-Q3 figures down to row 3 first.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.01

# Write the new 2022-Q4 summary figures into row 2.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# Step 2: turn the existing "2022-Q3" sheet into "2022-Q4", then duplicate it
# (before changing any data) so the original Q3 data/style survives unchanged
# on a new tab named "2022-Q3".
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

$wsQ4.Copy($null, $wsQ4)
$wsQ3 = $wb.Worksheets.Item(3)
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# Step 3: replace the data on the "2022-Q4" sheet with the new holdings.
# ---------------------------------------------------------------------------

# Clear out the previous Q3 holdings rows (keep header row 1 untouched).
$wsQ4.Rows.Item(2).Resize(2).ClearContents()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")

# Re-assert the header text (unchanged, but make sure it's present) and make
# sure the header row keeps the same style used elsewhere in the workbook.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ4.Range($cols[$i] + "1").Value = $headers[$i]
}
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)

# Data rows: A=index, B..G are kept as literal text (matches source data),
# H is numeric.
$rows = @(
    @(0, "005167", "嘉实润泽量化一年定期开放混合", "0.56", "27.25", "0.59", "0.0033", 7),
    @(1, "000926", "中信建投睿信灵活配置混合A", "0.10", "83.25", "1.04", "0.0010", 6),
    @(2, "004676", "中信建投睿信灵活配置混合C", "0.03", "83.25", "1.04", "0.0003", 6)
)

# Force columns B..G to be treated as plain text so values such as "0.56" or
# "005167" are not silently reinterpreted as numbers.
$wsQ4.Range("B2:G4").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $wsQ4.Range("A" + $excelRow).Value = $row[0]
    $wsQ4.Range("B" + $excelRow).Value = $row[1]
    $wsQ4.Range("C" + $excelRow).Value = $row[2]
    $wsQ4.Range("D" + $excelRow).Value = $row[3]
    $wsQ4.Range("E" + $excelRow).Value = $row[4]
    $wsQ4.Range("F" + $excelRow).Value = $row[5]
    $wsQ4.Range("G" + $excelRow).Value = $row[6]
    $wsQ4.Range("H" + $excelRow).Value = $row[7]
}

# Column A keeps the same highlighted style used throughout the workbook.
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("A4").Value = 2

Write-Host "2022-Q4 sheet populated; workbook now has" $wb.Worksheets.Count "sheets"
